$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formats from the last existing data row (115) onto the two new rows ---
$ws.Range("A115:V115").Copy()
$ws.Range("A116:V116").PasteSpecial(-4122)
$ws.Range("A115:V115").Copy()
$ws.Range("A117:V117").PasteSpecial(-4122)

# --- Row 116: Ameliano vs Olimpia Asuncion ---
$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = "paraguay"
$ws.Cells.Item(116, 3).Value = "primera-division"

$ws.Range("D116").NumberFormat = "@"
$ws.Cells.Item(116, 4).Value = "2023"
$ws.Range("D116").Style = "Normal"

$ws.Cells.Item(116, 5).Value = 45241.9375
$ws.Cells.Item(116, 6).Value = "Ameliano"
$ws.Cells.Item(116, 7).Value = 1
$ws.Cells.Item(116, 8).Value = "Olimpia Asuncion"
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 3.14
$ws.Cells.Item(116, 11).Value = "08/11/2023 14:42"
$ws.Cells.Item(116, 12).Value = 3.26
$ws.Cells.Item(116, 13).Value = "11/11/2023 22:26"
$ws.Cells.Item(116, 14).Value = 3.52
$ws.Cells.Item(116, 15).Value = "08/11/2023 14:42"
$ws.Cells.Item(116, 16).Value = 3.53
$ws.Cells.Item(116, 17).Value = "11/11/2023 22:26"
$ws.Cells.Item(116, 18).Value = 2.27
$ws.Cells.Item(116, 19).Value = "08/11/2023 14:42"
$ws.Cells.Item(116, 20).Value = 2.25
$ws.Cells.Item(116, 21).Value = "11/11/2023 22:26"
$ws.Cells.Item(116, 22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/sportivo-ameliano-olimpia-asuncion/UPViu3ZR/"

# --- Row 117: General Caballero JLM vs Guairena ---
$ws.Cells.Item(117, 1).Value = 116
$ws.Cells.Item(117, 2).Value = "paraguay"
$ws.Cells.Item(117, 3).Value = "primera-division"

$ws.Range("D117").NumberFormat = "@"
$ws.Cells.Item(117, 4).Value = "2023"
$ws.Range("D117").Style = "Normal"

$ws.Cells.Item(117, 5).Value = 45242.02083333334
$ws.Cells.Item(117, 6).Value = "General Caballero JLM"
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = "Guairena"
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 2.14
$ws.Cells.Item(117, 11).Value = "08/11/2023 14:42"
$ws.Cells.Item(117, 12).Value = 2.3
$ws.Cells.Item(117, 13).Value = "12/11/2023 00:23"
$ws.Cells.Item(117, 14).Value = 3.3
$ws.Cells.Item(117, 15).Value = "08/11/2023 14:42"
$ws.Cells.Item(117, 16).Value = 3.07
$ws.Cells.Item(117, 17).Value = "12/11/2023 00:23"
$ws.Cells.Item(117, 18).Value = 3.67
$ws.Cells.Item(117, 19).Value = "08/11/2023 14:42"
$ws.Cells.Item(117, 20).Value = 3.67
$ws.Cells.Item(117, 21).Value = "12/11/2023 00:23"
$ws.Cells.Item(117, 22).Value = "https://www.betexplorer.com/football/paraguay/primera-division/general-caballero-jlm-guairena-fc/MHNJzNsq/"
